# Derek's Log - add a FRIDAY section (separator row + one AV Shutdown entry)
# at the bottom of the "Logs" sheet, mirroring the existing WEDNESDAY section.
# (Off-by-one bug fix: the logged range now extends two rows further.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$yellow = 65535
$xlPasteFormats = -4122

# ---- Row 41: "FRIDAY" separator row ----
# Pick up the per-column number-format/border pattern from a plain data row
# (row 4) and just re-tint it yellow, same as the existing WEDNESDAY banner.
for ($col = 1; $col -le 6; $col++) {
    $ws.Cells.Item(4, $col).Copy() | Out-Null
    $ws.Cells.Item(41, $col).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Cells.Item(41, $col).Interior.Color = $yellow
}
$ws.Cells.Item(41, 3).Value = "FRIDAY"

# ---- Row 42: new AV Shutdown entry, identical shape/format to row 40 ----
for ($col = 1; $col -le 6; $col++) {
    $ws.Cells.Item(40, $col).Copy() | Out-Null
    $ws.Cells.Item(42, $col).PasteSpecial($xlPasteFormats) | Out-Null
}
$ws.Rows.Item(42).RowHeight = 30

$ws.Cells.Item(42, 1).Value = "AV Shutdown"
$ws.Cells.Item(42, 2).Value = 42587
$ws.Cells.Item(42, 3).Value = "1600"
$ws.Cells.Item(42, 4).Value = "FC"
$ws.Cells.Item(42, 5).Value = "305 SCR"
$ws.Cells.Item(42, 6).Value = "Turn off PC and projector. Leave equipment in room. Lock room. Key for room in Founders 164 storeroom."

$excel.CutCopyMode = $false

# ---- sheet bookkeeping: extend used range + move the last selection down ----
$ws.Range("F46").Select() | Out-Null
